# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 7106
$ws1.Range("F6").Value  = 556
$ws1.Range("F7").Value  = 162
$ws1.Range("F13").Value = 450
$ws1.Range("F17").Value = 3677
$ws1.Range("F21").Value = 27
$ws1.Range("F23").Value = 2315
$ws1.Range("F30").Value = 20
$ws1.Range("F32").Value = 1344
$ws1.Range("F33").Value = 121

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 7106
$ws4.Range("F7").Value  = 556
$ws4.Range("F8").Value  = 162
$ws4.Range("F14").Value = 450
$ws4.Range("F18").Value = 3677
$ws4.Range("F22").Value = 27
$ws4.Range("F24").Value = 2315
$ws4.Range("F31").Value = 20
$ws4.Range("F33").Value = 1344
$ws4.Range("F34").Value = 121
